# Fix clinical breakpoint seed
# - Wrong enum value for CLSI
# - Missing ValidFrom data
# - Missing "No-Specific" species entries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Breakpoints EUCAST and CLSI")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Missing "Non-species specific" -> "None" rename (Fluconazole rows for both EUCAST & CLSI)
$ws.Range("B44").Value = "None"
$ws.Range("B108").Value = "None"

# Missing ValidFrom data for the CLSI section (rows 66-129)
$validFromRange = $ws.Range("E66:E129")
$validFromRange.WrapText = $true
$validFromRange.HorizontalAlignment = -4108
$validFromRange.VerticalAlignment = -4108
$validFromRange.NumberFormat = "@"
$validFromRange.Value = "2019-11-01"

# Match the last selected cell recorded in the sheet view
$ws.Range("H9").Select()
